$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the data values (columns D, L, M, N, O, P, S) among rows 2-5:
#   new row2 = old row4, new row3 = old row5, new row4 = old row3, new row5 = old row2
# Capture the "before" values first, then write them back in the new order.

$rows = @(2, 3, 4, 5)
$cols = @("D", "L", "M", "N", "O", "P", "S")

$before = @{}
foreach ($r in $rows) {
    $before[$r] = @{}
    foreach ($c in $cols) {
        $before[$r][$c] = $ws.Range("$c$r").Value2
    }
}

$mapping = @{ 2 = 4; 3 = 5; 4 = 3; 5 = 2 }

foreach ($r in $rows) {
    $src = $mapping[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value2 = $before[$src][$c]
    }
}
